$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 89.25
$ws.Range("I2").Value = 72.28570999999999
$ws.Range("J2").Value = 113
$ws.Range("K2").Value = 72.28570999999999
$ws.Range("L2").Value = 113
$ws.Range("M2").Value = 40.71429000000001
$ws.Range("N2").Value = -339
$ws.Range("H12").Value = 188.29167
$ws.Range("J12").Value = 177.75
$ws.Range("L12").Value = 177.75
$ws.Range("N12").Value = -517.75
$ws.Range("H40").Value = 10526.333
$ws.Range("I40").Value = 25528
$ws.Range("J40").Value = 3025.5
$ws.Range("K40").Value = 25528
$ws.Range("L40").Value = 3025.5
$ws.Range("M40").Value = -25353
$ws.Range("N40").Value = -3375.5
$ws.Range("H64").Value = 54424.5
$ws.Range("I64").Value = 4500.3335
$ws.Range("K64").Value = 4500.3335
$ws.Range("M64").Value = -4252.3335
$ws.Range("H67").Value = 54424.5
$ws.Range("I67").Value = 4500.3335
$ws.Range("K67").Value = 4500.3335
$ws.Range("M67").Value = -3642.3335
$ws.Range("H76").Value = 22225488
$ws.Range("I76").Value = 25644096
$ws.Range("J76").Value = 4535.5
$ws.Range("K76").Value = 25644096
$ws.Range("L76").Value = 4535.5
$ws.Range("M76").Value = -25643781
$ws.Range("N76").Value = -5165.5
$ws.Range("H79").Value = 22225488
$ws.Range("I79").Value = 25644096
$ws.Range("J79").Value = 4535.5
$ws.Range("K79").Value = 25644096
$ws.Range("L79").Value = 4535.5
$ws.Range("M79").Value = -25643004
$ws.Range("N79").Value = -6719.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2968.625
$ws.Range("I61").Value = 2552
$ws.Range("J61").Value = 4033.3333
$ws.Range("K61").Value = 2552
$ws.Range("L61").Value = 4033.3333
$ws.Range("M61").Value = -2340
$ws.Range("N61").Value = -4457.3333
$ws.Range("H63").Value = 1974.8
$ws.Range("I63").Value = 1916.5555
$ws.Range("K63").Value = 1916.5555
$ws.Range("M63").Value = -1230.5555
$ws.Range("H66").Value = 1974.8
$ws.Range("I66").Value = 1916.5555
$ws.Range("K66").Value = 9582.7775
$ws.Range("M66").Value = -6150.7775
$ws.Range("H74").Value = 29596.805
$ws.Range("I74").Value = 45527.305
$ws.Range("J74").Value = 1412.0769
$ws.Range("K74").Value = 45527.305
$ws.Range("L74").Value = 1412.0769
$ws.Range("M74").Value = -44653.305
$ws.Range("N74").Value = -3160.0769
$ws.Range("H77").Value = 29596.805
$ws.Range("I77").Value = 45527.305
$ws.Range("J77").Value = 1412.0769
$ws.Range("K77").Value = 227636.525
$ws.Range("L77").Value = 7060.3845
$ws.Range("M77").Value = -223268.525
$ws.Range("N77").Value = -15796.3845
$ws.Range("H132").Value = 213055.64
$ws.Range("I132").Value = 34254
$ws.Range("J132").Value = 591694.4399999999
$ws.Range("K132").Value = 102762
$ws.Range("L132").Value = 1775083.32
$ws.Range("M132").Value = -100232
$ws.Range("N132").Value = -1780143.32
$ws.Range("H136").Value = 2968.625
$ws.Range("I136").Value = 2552
$ws.Range("J136").Value = 4033.3333
$ws.Range("K136").Value = 7656
$ws.Range("L136").Value = 12099.9999
$ws.Range("M136").Value = -5106
$ws.Range("N136").Value = -17199.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 14976.846
$ws.Range("I80").Value = 63338.668
$ws.Range("J80").Value = 468.3
$ws.Range("K80").Value = 63338.668
$ws.Range("L80").Value = 468.3
$ws.Range("M80").Value = -62340.668
$ws.Range("N80").Value = -2464.3
$ws.Range("H83").Value = 14976.846
$ws.Range("I83").Value = 63338.668
$ws.Range("J83").Value = 468.3
$ws.Range("K83").Value = 316693.34
$ws.Range("L83").Value = 2341.5
$ws.Range("M83").Value = -311701.34
$ws.Range("N83").Value = -12325.5
$ws.Range("H105").Value = 725202.6
$ws.Range("I105").Value = 1593065.8
$ws.Range("J105").Value = 1983.3334
$ws.Range("K105").Value = 1593065.8
$ws.Range("L105").Value = 1983.3334
$ws.Range("M105").Value = -1591318.8
$ws.Range("N105").Value = -5477.3334
$ws.Range("H134").Value = 25027060
$ws.Range("I134").Value = 1996.8572
$ws.Range("K134").Value = 5990.571599999999
$ws.Range("M134").Value = -3455.571599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 35721190
$ws.Range("I86").Value = 66676490
$ws.Range("J86").Value = 3541.6924
$ws.Range("K86").Value = 66676490
$ws.Range("L86").Value = 3541.6924
$ws.Range("M86").Value = -66675367
$ws.Range("N86").Value = -5787.6924
$ws.Range("H89").Value = 35721190
$ws.Range("I89").Value = 66676490
$ws.Range("J89").Value = 3541.6924
$ws.Range("K89").Value = 333382450
$ws.Range("L89").Value = 17708.462
$ws.Range("M89").Value = -333376834
$ws.Range("N89").Value = -28940.462
$ws.Range("H94").Value = 2619.5881
$ws.Range("I94").Value = 3496.5
$ws.Range("J94").Value = 2141.2727
$ws.Range("K94").Value = 3496.5
$ws.Range("L94").Value = 2141.2727
$ws.Range("M94").Value = -3045.5
$ws.Range("N94").Value = -3043.2727
$ws.Range("H122").Value = 1266.381
$ws.Range("I122").Value = 1495.6923
$ws.Range("J122").Value = 893.75
$ws.Range("K122").Value = 4487.0769
$ws.Range("L122").Value = 2681.25
$ws.Range("M122").Value = -2037.0769
$ws.Range("N122").Value = -7581.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 414.58823
$ws.Range("I113").Value = 378.66666
$ws.Range("J113").Value = 500.8
$ws.Range("K113").Value = 1135.99998
$ws.Range("L113").Value = 1502.4
$ws.Range("M113").Value = 1034.00002
$ws.Range("N113").Value = -5842.4
$ws.Range("H131").Value = 17596038
$ws.Range("J131").Value = 21391196
$ws.Range("L131").Value = 64173588
$ws.Range("N131").Value = -64183668
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5847.9
$ws.Range("I80").Value = 2729.3076
$ws.Range("J80").Value = 8232.706
$ws.Range("K80").Value = 2729.3076
$ws.Range("L80").Value = 8232.706
$ws.Range("M80").Value = -1731.3076
$ws.Range("N80").Value = -10228.706
$ws.Range("H83").Value = 5847.9
$ws.Range("I83").Value = 2729.3076
$ws.Range("J83").Value = 8232.706
$ws.Range("K83").Value = 13646.538
$ws.Range("L83").Value = 41163.53
$ws.Range("M83").Value = -8654.538
$ws.Range("N83").Value = -51147.53
$ws.Range("H140").Value = 60535.8
$ws.Range("J140").Value = 60535.8
$ws.Range("L140").Value = 60535.8
$ws.Range("N140").Value = -70895.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1038.3846
$ws.Range("I16").Value = 833.875
$ws.Range("J16").Value = 1365.6
$ws.Range("K16").Value = 833.875
$ws.Range("L16").Value = 1365.6
$ws.Range("M16").Value = -663.875
$ws.Range("N16").Value = -1705.6
$ws.Range("H68").Value = 2255
$ws.Range("I68").Value = 1503
$ws.Range("J68").Value = 2589.2222
$ws.Range("K68").Value = 1503
$ws.Range("L68").Value = 2589.2222
$ws.Range("M68").Value = -754
$ws.Range("N68").Value = -4087.2222
$ws.Range("H71").Value = 2255
$ws.Range("I71").Value = 1503
$ws.Range("J71").Value = 2589.2222
$ws.Range("K71").Value = 7515
$ws.Range("L71").Value = 12946.111
$ws.Range("M71").Value = -3771
$ws.Range("N71").Value = -20434.111
$ws.Range("H136").Value = 295365.97
$ws.Range("I136").Value = 417513.25
$ws.Range("J136").Value = 2212.5
$ws.Range("K136").Value = 1252539.75
$ws.Range("L136").Value = 6637.5
$ws.Range("M136").Value = -1249989.75
$ws.Range("N136").Value = -11737.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 25000
$ws.Range("I94").Value = 25000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 25000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -24099
$ws.Range("N94").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 48000
$ws.Range("J111").Value = 48000
$ws.Range("L111").Value = 48000
$ws.Range("N111").Value = -56180
$ws.Range("H124").Value = 16000
$ws.Range("J124").Value = 16000
$ws.Range("L124").Value = 16000
$ws.Range("N124").Value = -25820
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H140").Value = 39871.4
$ws.Range("J140").Value = 39871.4
$ws.Range("L140").Value = 39871.4
$ws.Range("N140").Value = -50231.4
